$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '27.247.12'
$ws.Range('D2').ClearFormats()
$ws.Range('E2').Value = '  +1.55%  '
$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '1.644.36'
$ws.Range('D3').ClearFormats()
$ws.Range('E3').Value = '  +0.26%  '
$ws.Range('E4').Value = '  +0.16%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '217.33'
$ws.Range('D5').ClearFormats()
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '0.518'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  +1.02%  '
$ws.Range('E7').Value = '  +0.18%  '
$ws.Range('E8').Value = '  +0.02%  '
$ws.Range('E9').Value = '  +0.97%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '20.03'
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = '  +0.60%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.0849'
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  +0.37%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '1.875.74'
$ws.Range('D12').ClearFormats()
$ws.Range('E12').Value = '  +0.39%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '1.643.94'
$ws.Range('D13').ClearFormats()
$ws.Range('E13').Value = '  +0.22%  '
$ws.Range('E14').Value = '  +0.81%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '0.544'
$ws.Range('D15').ClearFormats()
$ws.Range('E15').Value = '  +2.84%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '67.02'
$ws.Range('D16').ClearFormats()
$ws.Range('E16').Value = '  +0.44%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '27.249.45'
$ws.Range('D17').ClearFormats()
$ws.Range('E17').Value = '  +1.56%  '
$ws.Range('E18').Value = '  +1.67%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '220.00'
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '  +0.33%  '
$ws.Range('E20').Value = '  +0.13%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '6.97'
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  +3.71%  '
$ws.Range('E22').Value = '  +3.84%  '
$ws.Range('E23').Value = '  +0.42%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '9.14'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  -0.29%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '148.84'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  +1.02%  '
$ws.Range('B26').Value = 'BinanceUSD'
$ws.Range('C26').Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '1.00'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  +0.19%  '
$ws.Range('B27').Value = 'Cosmos'
$ws.Range('C27').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '7.49'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  +1.64%  '
$ws.Range('E28').Value = '  -0.48%  '
$ws.Range('E29').Value = '  -0.16%  '
$ws.Range('E30').Value = '  +1.80%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '1.19'
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  +0.33%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '3.39'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  +1.77%  '
$ws.Range('E33').Value = '  +0.36%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '1.307.40'
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  +3.74%  '
$ws.Range('E35').Value = '  +1.33%  '
$ws.Range('E36').Value = '  +1.63%  '
$ws.Range('E37').Value = '  -0.31%  '
$ws.Range('E38').Value = '  +3.70%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '0.860'
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  +3.28%  '
$ws.Range('E40').Value = '  +0.11%  '
$ws.Range('E41').Value = '  +0.55%  '
$ws.Range('E42').Value = '  +5.76%  '
$ws.Range('E43').Value = '  -2.47%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '1.784.76'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  +0.20%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '61.95'
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  +0.25%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '91.96'
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '  +0.08%  '
$ws.Range('E47').Value = '  +1.82%  '
$ws.Range('E48').Value = '  +2.05%  '
$ws.Range('E49').Value = '  -0.06%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '7.71'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  +0.85%  '
$ws.Range('E51').Value = '  +0.34%  '
